## "hybrid inout EXP version"
##
## - Tweaks the view of the existing "inflow vs outflow" sheet (zoom 85%,
##   selection C3:E57, no longer the selected tab).
## - Adds a new "hybrid InOut" sheet after it, with the same C:E mini-table
##   skeleton (row numbers 1..50, avg/SD formulas) but with the train/test
##   columns still empty (so avg/SD read as #DIV/0!), and makes it the
##   active/selected sheet (zoom 70%, selection L31).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---- 1. Update the view of the original sheet while it is still active ----
$ws1.Activate()
$excel.ActiveWindow.Zoom = 85
$ws1.Range("C3:E57").Select()

# ---- 2. Add the new sheet right after "inflow vs outflow" ----
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "hybrid InOut"

# Column widths for D (train) and E (test)
$ws2.Range("D1").ColumnWidth = 11.76
$ws2.Range("E1").ColumnWidth = 13.25

# ---- 3. Build the C3:E57 skeleton, re-using sheet1's look (fonts/number format) ----
$ws1.Range("C3").Copy()
$ws2.Range("C3").PasteSpecial(-4122)
$ws1.Range("E3").Copy()
$ws2.Range("E3").PasteSpecial(-4122)

$ws1.Range("C4:E4").Copy()
$ws2.Range("C4").PasteSpecial(-4122)

$ws1.Range("C5:E54").Copy()
$ws2.Range("C5").PasteSpecial(-4122)

$ws1.Range("C56:E57").Copy()
$ws2.Range("C56").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Header / label text
$ws2.Range("C3").Value = "Random Forest-100 (superdataset-24 InOut.csv)"
$ws2.Range("D4").Value = "train (MAE)"
$ws2.Range("E4").Value = "test (MAE)"
$ws2.Range("C56").Value = "avg"
$ws2.Range("C57").Value = "SD"

# Row counter 1..50 down column C
$ws2.Range("C5").Value = 1
$ws2.Range("C6").Formula = "=C5+1"
$ws2.Range("C7:C54").FormulaR1C1 = "=R[-1]C+1"

# avg / SD formulas for the (still empty) train/test columns
$ws2.Range("D56").Formula = "=AVERAGE(D5:D54)"
$ws2.Range("E56").Formula = "=AVERAGE(E5:E54)"
$ws2.Range("D57").Formula = "=STDEV.S(D5:D54)"
$ws2.Range("E57").Formula = "=STDEV.S(E5:E54)"

# ---- 4. View state for the new sheet: it becomes the active/selected tab ----
$excel.ActiveWindow.Zoom = 70
$ws2.Range("L31").Select()
